# Atualização de bases das ligas, do dia: 08-04-2024 às 21:28
#
# The underlying data rows got re-sorted (pairs of adjacent match rows swapped
# places), which is why entire rows B:AC (everything except the running id in
# column A) need to trade places. Using Value2 both to read the whole row
# (as a 2D array) and to write it back keeps all data types (numbers, dates,
# strings) intact in one shot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AC$r1")
    $rng2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

Swap-Rows 3 4
Swap-Rows 5 7
Swap-Rows 16 17
Swap-Rows 20 21
Swap-Rows 46 47
Swap-Rows 67 68
Swap-Rows 86 87
Swap-Rows 120 121
